# Commit: "added new trials for submission with more covid data"
# Adds 13 new submission-trial rows (sheet rows 33-45) to the Tabelle2
# table on worksheet "Tabelle1", growing it from A1:E33 to A1:E45.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$table = $ws.ListObjects.Item(1)

# Columns: Date(serial), Name, Name Ramp, Hand in, By
$newRows = @(
    @(44894, "221129_final_model",     "final_model",          "TRUE", "Maria"),
    @(44894, "221129_final_model_v2",  "final_model_v2",       "TRUE", "Maria"),
    @(44894, "221129_test_new_data",   "test_new_data",        "TRUE", "Maria"),
    @(44894, "221129_test_new_data_2", "model_name",           "TRUE", "Maria"),
    @(44894, "-",                      "model_name_2",         "TRUE", "Maria"),
    @(44894, "221129_test_new_data_3", "14_gewinnt",           "TRUE", "Maria"),
    @(44894, "221129_test_new_data_4", "caca_aos_gambuzinos",  "TRUE", "Maria"),
    @(44894, "221129_test_new_data_5", "estimator_1",          "TRUE", "Maria"),
    @(44894, "221129_test_new_data_6", "estimator_2",          "TRUE", "Maria"),
    @(44894, "-",                      "estimator_3",          "TRUE", "Maria"),
    @(44894, "221129_test_new_data_7", "estimator_4",          "TRUE", "Maria"),
    @(44895, "221130_best_grid",       "estimator_5",          "TRUE", "Maria"),
    @(44895, "221130_lowest_test",     "estimaror_6",          "TRUE", "Maria")
)

$startRow = 33
$dateFormatSrc = $ws.Cells.Item($startRow - 1, 1)   # an existing date cell to clone number format from
$trueCellSrc   = $ws.Cells.Item($startRow - 1, 4)   # an existing "TRUE" cell (text, not boolean)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Row 33 already exists as the (blank) last row of the table; grow the
    # table by one row for every row after that so the ref/autoFilter range
    # tracks along correctly (ends up at A1:E45, not A1:E46).
    if ($r -gt $startRow) {
        $table.ListRows.Add() | Out-Null
    }

    # Column A: date, formatted like the other date cells in the column.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $data[0]
    $dateFormatSrc.Copy()
    $cellA.PasteSpecial(-4122)  # xlPasteFormats

    # Column B: Name
    $ws.Cells.Item($r, 2).Value = $data[1]

    # Column C: Name Ramp
    $ws.Cells.Item($r, 3).Value = $data[2]

    # Column D: Hand in - stored as literal text "TRUE" (shared string),
    # not a real boolean, in this sheet. Copy the value from an existing
    # "TRUE" text cell so it lands as text without forcing a quote-prefix
    # cell style.
    $trueCellSrc.Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4163)  # xlPasteValues

    # Column E: By
    $ws.Cells.Item($r, 5).Value = $data[4]
}

$excel.CutCopyMode = $false

$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.Zoom = 100
$ws.Range("B44").Select() | Out-Null
